$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 32211.844
$ws.Range("J17").Value = 32211.844
$ws.Range("L17").Value = 96635.53200000001
$ws.Range("N17").Value = -96971.53200000001

$ws.Range("H29").Value = 3332.111
$ws.Range("I29").Value = 96.333336
$ws.Range("J29").Value = 4950
$ws.Range("K29").Value = 289.000008
$ws.Range("L29").Value = 14850
$ws.Range("M29").Value = -8.00000799999998
$ws.Range("N29").Value = -15412

$ws.Range("H38").Value = 981.58826
$ws.Range("I38").Value = 71.545456
$ws.Range("J38").Value = 2650
$ws.Range("K38").Value = 214.636368
$ws.Range("L38").Value = 7950
$ws.Range("M38").Value = 157.363632
$ws.Range("N38").Value = -8694

$ws.Range("H58").Value = 2382
$ws.Range("I58").Value = 707.5
$ws.Range("J58").Value = 2621.2144
$ws.Range("K58").Value = 2122.5
$ws.Range("L58").Value = 7863.6432
$ws.Range("M58").Value = -1972.5
$ws.Range("N58").Value = -8163.6432

$ws.Range("H62").Value = 2239.7693
$ws.Range("I62").Value = 1645
$ws.Range("J62").Value = 2611.5
$ws.Range("K62").Value = 1645
$ws.Range("L62").Value = 2611.5
$ws.Range("M62").Value = -1021
$ws.Range("N62").Value = -3859.5

$ws.Range("H65").Value = 2239.7693
$ws.Range("I65").Value = 1645
$ws.Range("J65").Value = 2611.5
$ws.Range("K65").Value = 8225
$ws.Range("L65").Value = 13057.5
$ws.Range("M65").Value = -5105
$ws.Range("N65").Value = -19297.5

$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").Value = ""

$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").Value = ""

$ws.Range("H87").Value = 25672.223
$ws.Range("J87").Value = 25672.223
$ws.Range("L87").Value = 25672.223
$ws.Range("N87").Value = -28168.223

$ws.Range("H90").Value = 25672.223
$ws.Range("J90").Value = 25672.223
$ws.Range("L90").Value = 77016.66900000001
$ws.Range("N90").Value = -89496.66900000001

$ws.Range("H98").Value = 3429.6758
$ws.Range("I98").Value = 2102.8
$ws.Range("J98").Value = 26650
$ws.Range("K98").Value = 2102.8
$ws.Range("L98").Value = 26650
$ws.Range("M98").Value = -604.8000000000002
$ws.Range("N98").Value = -29646

$ws.Range("H122").Value = 3429.6758
$ws.Range("I122").Value = 2102.8
$ws.Range("J122").Value = 26650
$ws.Range("K122").Value = 6308.400000000001
$ws.Range("L122").Value = 79950
$ws.Range("M122").Value = -3858.400000000001
$ws.Range("N122").Value = -84850

$ws.Range("H129").Value = 1013.4681
$ws.Range("J129").Value = 1095.0769
$ws.Range("L129").Value = 3285.2307
$ws.Range("N129").Value = -13285.2307

$ws.Range("H138").Value = 4203.5713
$ws.Range("I138").Value = 2203.3022
$ws.Range("J138").Value = 8504.15
$ws.Range("K138").Value = 6609.9066
$ws.Range("L138").Value = 25512.45
$ws.Range("M138").Value = -1469.9066
$ws.Range("N138").Value = -35792.45

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12519.846
$ws.Range("I32").Value = 15039.415
$ws.Range("K32").Value = 15039.415
$ws.Range("M32").Value = -14752.415

$ws.Range("H61").Value = 1520.091
$ws.Range("I61").Value = 1346.2174
$ws.Range("J61").Value = 1920
$ws.Range("K61").Value = 1346.2174
$ws.Range("L61").Value = 1920
$ws.Range("M61").Value = -1134.2174
$ws.Range("N61").Value = -2344

$ws.Range("H97").Value = 685.3913
$ws.Range("I97").Value = 664.6875
$ws.Range("J97").Value = 732.7143
$ws.Range("K97").Value = 664.6875
$ws.Range("L97").Value = 732.7143
$ws.Range("M97").Value = -168.6875
$ws.Range("N97").Value = -1724.7143

$ws.Range("H110").Value = 905.3158
$ws.Range("I110").Value = 967.1724
$ws.Range("J110").Value = 706
$ws.Range("K110").Value = 967.1724
$ws.Range("L110").Value = 706
$ws.Range("M110").Value = 1077.8276
$ws.Range("N110").Value = -4796

$ws.Range("H121").Value = 30418.334
$ws.Range("J121").Value = 30418.334
$ws.Range("L121").Value = 30418.334
$ws.Range("N121").Value = -33912.334

$ws.Range("H132").Value = 1693.2373
$ws.Range("I132").Value = 1043.8
$ws.Range("K132").Value = 3131.4
$ws.Range("M132").Value = -601.3999999999996

$ws.Range("H133").Value = 73254.86
$ws.Range("J133").Value = 73254.86
$ws.Range("L133").Value = 73254.86
$ws.Range("N133").Value = -78314.86

$ws.Range("H136").Value = 1520.091
$ws.Range("I136").Value = 1346.2174
$ws.Range("J136").Value = 1920
$ws.Range("K136").Value = 4038.6522
$ws.Range("L136").Value = 5760
$ws.Range("M136").Value = -1488.6522
$ws.Range("N136").Value = -10860

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 21264
$ws.Range("I5").Value = 6996
$ws.Range("J5").Value = 49800
$ws.Range("K5").Value = 6996
$ws.Range("L5").Value = 49800
$ws.Range("M5").Value = -6883
$ws.Range("N5").Value = -50026

$ws.Range("H94").Value = 85128.086
$ws.Range("I94").Value = 943.375
$ws.Range("J94").Value = 253497.5
$ws.Range("K94").Value = 943.375
$ws.Range("L94").Value = 253497.5
$ws.Range("M94").Value = -492.375
$ws.Range("N94").Value = -254399.5

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = ""

$ws.Range("H107").Value = 25917.238
$ws.Range("J107").Value = 2952
$ws.Range("L107").Value = 2952
$ws.Range("N107").Value = -6792

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 295081.7
$ws.Range("I132").Value = 398472.38
$ws.Range("J132").Value = 2141.4167
$ws.Range("K132").Value = 1195417.14
$ws.Range("L132").Value = 6424.250100000001
$ws.Range("M132").Value = -1192887.14
$ws.Range("N132").Value = -11484.2501

$ws.Range("H135").Value = 71000
$ws.Range("J135").Value = 71000
$ws.Range("L135").Value = 71000
$ws.Range("N135").Value = -81140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 250
$ws.Range("I4").Value = 250
$ws.Range("K4").Value = 750
$ws.Range("M4").Value = -638

$ws.Range("H113").Value = 909.875
$ws.Range("I113").Value = 1001.5
$ws.Range("J113").Value = 879.3333
$ws.Range("K113").Value = 3004.5
$ws.Range("L113").Value = 2637.9999
$ws.Range("M113").Value = -834.5
$ws.Range("N113").Value = -6977.9999

$ws.Range("H131").Value = 13172617
$ws.Range("I131").Value = 15961.429
$ws.Range("J131").Value = 14507350
$ws.Range("K131").Value = 47884.287
$ws.Range("L131").Value = 43522050
$ws.Range("M131").Value = -42844.287
$ws.Range("N131").Value = -43532130

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = ""

$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").Value = ""

$ws.Range("H63").Value = 30466.666
$ws.Range("J63").Value = 30466.666
$ws.Range("L63").Value = 30466.666
$ws.Range("N63").Value = -31838.666

$ws.Range("H66").Value = 30466.666
$ws.Range("J66").Value = 30466.666
$ws.Range("L66").Value = 91399.99800000001
$ws.Range("N66").Value = -98263.99800000001

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").Value = ""

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").Value = ""

$ws.Range("H97").Value = 116978.46
$ws.Range("I97").Value = 43393.332
$ws.Range("J97").Value = 1000000
$ws.Range("K97").Value = 43393.332
$ws.Range("L97").Value = 1000000
$ws.Range("M97").Value = -42897.332
$ws.Range("N97").Value = -1000992

$ws.Range("H122").Value = 4150
$ws.Range("I122").Value = 3500
$ws.Range("J122").Value = 4800
$ws.Range("K122").Value = 10500
$ws.Range("L122").Value = 14400
$ws.Range("M122").Value = -8050
$ws.Range("N122").Value = -19300

$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = ""
$ws.Range("N126").Value = -19940

$ws.Range("H132").Value = 1141.425
$ws.Range("I132").Value = 734.5
$ws.Range("J132").Value = 2090.9167
$ws.Range("K132").Value = 2203.5
$ws.Range("L132").Value = 6272.750100000001
$ws.Range("M132").Value = 326.5
$ws.Range("N132").Value = -11332.7501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2522.7646
$ws.Range("I16").Value = 2578.0833
$ws.Range("J16").Value = 2390
$ws.Range("K16").Value = 2578.0833
$ws.Range("L16").Value = 2390
$ws.Range("M16").Value = -2408.0833
$ws.Range("N16").Value = -2730

$ws.Range("H55").Value = 280
$ws.Range("I55").Value = 100
$ws.Range("J55").Value = 352
$ws.Range("K55").Value = 100
$ws.Range("L55").Value = 352
$ws.Range("M55").Value = 73
$ws.Range("N55").Value = -698

$ws.Range("H93").Value = 867.6875
$ws.Range("I93").Value = 571.1818
$ws.Range("K93").Value = 571.1818
$ws.Range("M93").Value = 676.8182

$ws.Range("H122").Value = 12506023
$ws.Range("I122").Value = 10875001
$ws.Range("J122").Value = 15391677
$ws.Range("K122").Value = 32625003
$ws.Range("L122").Value = 46175031
$ws.Range("M122").Value = -32622553
$ws.Range("N122").Value = -46179931

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 23112.75
$ws.Range("J64").Value = 23112.75
$ws.Range("L64").Value = 23112.75
$ws.Range("N64").Value = -23608.75

$ws.Range("H67").Value = 23112.75
$ws.Range("J67").Value = 23112.75
$ws.Range("L67").Value = 23112.75
$ws.Range("N67").Value = -24828.75

$ws.Range("H107").Value = 706.44446
$ws.Range("I107").Value = 682.3
$ws.Range("J107").Value = 736.625
$ws.Range("K107").Value = 2046.9
$ws.Range("L107").Value = 2209.875
$ws.Range("M107").Value = -126.8999999999999
$ws.Range("N107").Value = -6049.875

$ws.Range("H123").Value = 22292.625
$ws.Range("J123").Value = 22292.625
$ws.Range("L123").Value = 22292.625
$ws.Range("N123").Value = -32092.625
